# Auto-generated script to apply cryptos price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.101.42"
$ws.Range("E2").Value = "  -1.78%  "

$ws.Range("D3").Value = "2.437.24"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.74"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.30"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -2.49%  "

$ws.Range("D9").Value = "2.435.07"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("E10").Value = "  -6.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  -1.80%  "

$ws.Range("E12").Value = "  -5.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.74"
$ws.Range("E13").Value = "  -3.14%  "

$ws.Range("D14").Value = "2.884.43"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").Value = "68.033.37"
$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("E16").Value = "  -4.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "22.95"
$ws.Range("E17").Value = "  -4.91%  "

$ws.Range("D18").Value = "2.434.04"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("E19").Value = "  -3.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.73"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.05"
$ws.Range("E21").Value = "  -4.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.71"
$ws.Range("E22").Value = "  -3.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  -4.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.24"
$ws.Range("E25").Value = "  -4.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.61"
$ws.Range("E27").Value = "  -6.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  -7.16%  "

$ws.Range("D30").Value = "0.0₃0816"
$ws.Range("E30").Value = "  -5.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  -8.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "420.77"
$ws.Range("E33").Value = "  -4.59%  "

$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("E35").Value = "  -5.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.01"
$ws.Range("E36").Value = "  +0.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.99"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("E39").Value = "  -4.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.68"
$ws.Range("E40").Value = "  -2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.300"
$ws.Range("E41").Value = "  -3.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.31"
$ws.Range("E42").Value = "  -5.30%  "

$ws.Range("E43").Value = "  -6.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.07"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.03"
$ws.Range("E45").Value = "  -5.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.91"
$ws.Range("E46").Value = "  -4.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  -3.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0712"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.474"
$ws.Range("E49").Value = "  -7.06%  "

$ws.Range("E50").Value = "  -2.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0901"
$ws.Range("E51").Value = "  -1.81%  "

